$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.35
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 3.25
$ws.Range("L2").Value = 4.33
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("X2").Value = 9.5
$ws.Range("Z2").Value = 23
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 15
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 41
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 15
$ws.Range("AS2").Value = 351
$ws.Range("BA2").Value = 151
$ws.Range("BD2").Value = 151

# Row 3 updates
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65

# Row 4 updates
$ws.Range("G4").Value = 1.95
$ws.Range("I4").Value = 4.5
$ws.Range("AA4").Value = 19
$ws.Range("AE4").Value = 19
